$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daz Material Types")

# Insert a new column before column C ("Daz Material Method Prefix"),
# shifting old C -> D and old D -> E. Formulas referencing column C
# automatically re-point at column D.
$ws.Columns.Item(3).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 3).Value = "Comments"

# Fill in the new comments. Order matters here so that the underlying
# shared-string table is populated in the same sequence as the source
# workbook (DiffuseStrength row first, then ColorMap, then DiffuseColor).
$ws.Cells.Item(21, 3).Value = "Percent multiplied with diffuse colour (from map or DiffuseColor property)"
$ws.Cells.Item(5, 3).Value = "Diffuse Colour map"
$ws.Cells.Item(6, 3).Value = "Diffuse Colour is multiplied with Diffuse Map (if present) and Diffuse Strength"

# Match the italic "comment" styling used for the new column.
$c2 = $ws.Cells.Item(2, 3)
$c2.Font.Bold = $false
$c2.Font.Italic = $true
$c2.Font.Size = 10

$bodyRange = $ws.Range($ws.Cells.Item(3, 3), $ws.Cells.Item(48, 3))
$bodyRange.Font.Bold = $false
$bodyRange.Font.Italic = $true
$bodyRange.Font.Size = 10

# Approximate the widened column C (best achievable given engine rounding).
$ws.Columns.Item(3).ColumnWidth = 60.2

# Move the selection/cursor as in the source edit.
$ws.Range("A21").Select()

$wb.Save()
